$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B2 to "7C#14" (was "2C#2")
$ws.Range("B2").Value = "7C#14"

# Add new row 3 duplicating the original row 2 content, with Pull # = 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2C#2"
$ws.Range("C3").Value = "EXPRESS"
$ws.Range("D3").Value = "100+00"
$ws.Range("E3").Value = "200+00"

# Update selection to E5 (per diff)
$ws.Range("E5").Select()
